# Add "Item category" column to the import header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before J (shifts J:W -> K:X) and label it "Item category".
$ws.Columns("J:J").Insert()
$ws.Range("J1").Value = "Item category"

# Remove the two "|" placeholder columns (now at L and, after that removal, W).
$ws.Columns("L:L").Delete()
$ws.Columns("V:V").Delete()

$ws.Range("I1").Select()
